$wb = $excel.ActiveWorkbook

# --- CardData (sheet2): add row 9 with ID "8" ---
$wsCard = $wb.Worksheets.Item("CardData")
$wsCard.Range("A9").Value = "8"

# --- NameData (sheet1): add row 7 with ID "6" and name "Carlos Jacinta" ---
$wsName = $wb.Worksheets.Item("NameData")
$wsName.Range("A7").Value = "6"
$wsName.Range("E7").Value = "Carlos Jacinta"

# --- UDFData (sheet6): add row 5, a modified copy of row 3 (ID "2" -> "4") ---
$wsUdf = $wb.Worksheets.Item("UDFData")
$wsUdf.Range("A5").Value = "4"
$wsUdf.Range("B5:K5").Style = "Normal"
$wsUdf.Range("B5").Value = "udf data 1 modified"
$wsUdf.Range("C5").Value = "udf data 2"
$wsUdf.Range("D5").Value = "udf data 3"
$wsUdf.Range("E5").Value = "Sweet"
$wsUdf.Range("F5").Value = "Sour"
$wsUdf.Range("G5").Value = "udf data 6"
$wsUdf.Range("H5").Value = "udf data 7"
$wsUdf.Range("I5").Value = "udf data 8"
$wsUdf.Range("J5").Value = "udf data 9"
$wsUdf.Range("K5").Value = "udf data 10"

# --- ACHData (sheet3): add row 8 with ID "7" ---
$wsAch = $wb.Worksheets.Item("ACHData")
$wsAch.Range("A8").Value = "7"

# --- Selections per sheet (restores the cursor position saved in each sheet) ---
$wsCard.Range("A9").Select() | Out-Null
$wsName.Range("E7").Select() | Out-Null
$wsAch.Range("B7").Select() | Out-Null

# UDFData ends up as the active/selected sheet+cell, matching activeTab="5" and tabSelected on sheet6
$wsUdf.Range("B5").Select() | Out-Null
